{"js": "// Find the \"Edison Achalma\" paragraph that uses the \"Author\" style, then\n// insert a new Author-styled paragraph right after it containing the\n// institutional affiliation line.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.style === \"Author\" && p.text.trim() === \"Edison Achalma\") {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Edison Achalma' Author paragraph\");\n}\n\nconst affiliation =\n  \"Escuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\";\n\n// Build a minimal OOXML package containing a single Author-styled paragraph\n// and insert it right after the target paragraph. Using insertOoxml (rather\n// than insertParagraph/insertText) keeps the existing \"Edison Achalma\" run\n// untouched and produces a plain <w:t xml:space=\"preserve\"> run, matching\n// the document's existing paragraph markup exactly.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr><w:pStyle w:val=\"Author\"/></w:pPr>\n            <w:r><w:t xml:space=\"preserve\">${affiliation}</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst insertionRange = target.getRange(Word.RangeLocation.end);\ninsertionRange.insertOoxml(ooxml, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$affiliation = \"Escuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\"\n\n# Find the \"Edison Achalma\" paragraph that uses the \"Author\" style.\n$target = $null\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Style.NameLocal -eq \"Author\" -and $p.Range.Text.Trim() -eq \"Edison Achalma\") {\n        $target = $p\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'Edison Achalma' Author paragraph\"\n}\n\n# Insert a new paragraph right after it (same \"Author\" style) holding the\n# institutional affiliation line. Appending \"<CR>text\" to the (non-\n# collapsed) paragraph Range splits the paragraph mark so the existing\n# \"Edison Achalma\" run is left untouched and the new text becomes its own\n# paragraph that inherits the Author style.\n[void]$target.Range.InsertAfter([char]13 + $affiliation)\n\n# Re-fetch the freshly created paragraph and normalize its run markup\n# (xml:space=\"preserve\") to match the rest of the document's plain-text\n# runs by round-tripping it through WordOpenXML, scoped strictly to that\n# paragraph's own range so neighboring paragraphs are untouched.\n$newPara = $d.Paragraphs.Item($targetIndex + 1)\n$newRange = $newPara.Range\n$ooxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr><w:pStyle w:val=\"Author\"/></w:pPr>\n            <w:r><w:t xml:space=\"preserve\">$affiliation</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n[void]$newRange.InsertXML($ooxml)\n"}
